# Insert a new data row at row 470 (pushing the existing rows 470-566 down
# to 471-567). The new row 470 is a copy of the (old) row 470's data, but
# with a new "Fecha" (date serial) value of 45275.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(470).Insert()

$ws.Cells.Item(470, 1).Value  = 3
$ws.Cells.Item(470, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(470, 3).Value  = "Coquimbo"
$ws.Cells.Item(470, 4).Value  = 45275
$ws.Cells.Item(470, 5).Value  = 5
$ws.Cells.Item(470, 6).Value  = 100112001
$ws.Cells.Item(470, 7).Value  = "Berenjena"
$ws.Cells.Item(470, 8).Value  = "Sin especificar"
$ws.Cells.Item(470, 9).Value  = "Primera"
$ws.Cells.Item(470, 10).Value = 60
$ws.Cells.Item(470, 11).Value = 9000
$ws.Cells.Item(470, 12).Value = 9000
$ws.Cells.Item(470, 13).Value = 9000
$ws.Cells.Item(470, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(470, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(470, 16).Value = 150
$ws.Cells.Item(470, 17).Value = 60
$ws.Cells.Item(470, 18).Value = "Hortaliza"
